$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Ali bin Ahmad"
$ws.Range("B2").Value = "qD9PBv5QNbVNcR"
$ws.Range("E2").Value = "012-3456789"
$ws.Range("A3").Value = "Siti binti Aminah"
$ws.Range("B3").Value = "fetLua2MOOClsl"
$ws.Range("E3").Value = "013-4567890"
$ws.Range("A4").Value = "Lim Wei Chong"
$ws.Range("B4").Value = "rtmyXBHhwbnAUG"
$ws.Range("E4").Value = "014-5678901"
$ws.Range("A5").Value = "Preeti Kaur"
$ws.Range("B5").Value = "QzJaY2L4EyW6Pj"
$ws.Range("E5").Value = "015-6789012"
$ws.Range("A6").Value = "Muhammad Hafiz"
$ws.Range("B6").Value = "S6yetJh9pxRYWj"
$ws.Range("E6").Value = "016-7890123"
$ws.Range("A7").Value = "Chong Mei Ling"
$ws.Range("B7").Value = "egXGWEqM5w1awC"
$ws.Range("E7").Value = "017-8901234"
$ws.Range("A8").Value = "Arun Raj"
$ws.Range("B8").Value = "l1wLSXU32WIjUX"
$ws.Range("E8").Value = "018-9012345"
$ws.Range("A9").Value = "Farah binti Zain"
$ws.Range("B9").Value = "M3hvlWW3m3HFAQ"
$ws.Range("E9").Value = "019-0123456"
$ws.Range("A10").Value = "Tan Ah Heng"
$ws.Range("B10").Value = "BTdhIhi0T9f2Pe"
$ws.Range("E10").Value = "010-2345678"
$ws.Range("A11").Value = "Anusha Devi"
$ws.Range("B11").Value = "eLGzmdMyIrQGkD"
$ws.Range("E11").Value = "011-3456789"
$ws.Range("A12").Value = "Nur Aisyah"
$ws.Range("B12").Value = "1KGdZK7P5UsjsI"
$ws.Range("E12").Value = "012-5678901"
$ws.Range("A13").Value = "Wong Siew Mei"
$ws.Range("B13").Value = "gfAVDdxUVzAJME"
$ws.Range("E13").Value = "013-6789012"
$ws.Range("A14").Value = "Thiru Selvan"
$ws.Range("B14").Value = "rbjLa11L9LiYXn"
$ws.Range("E14").Value = "014-7890123"
$ws.Range("A15").Value = "Aina binti Razak"
$ws.Range("B15").Value = "Jb4pL17VqHUYuX"
$ws.Range("E15").Value = "015-8901234"
$ws.Range("A16").Value = "Yong Kai Wen"
$ws.Range("B16").Value = "nEbs5LCkAUFVFs"
$ws.Range("E16").Value = "016-9012345"
$ws.Range("A17").Value = "Hema Latha"
$ws.Range("B17").Value = "i10QF4IzAoCK9V"
$ws.Range("E17").Value = "017-0123456"
$ws.Range("A18").Value = "Ahmad bin Abu"
$ws.Range("B18").Value = "HOLEs7GWMcEPsr"
$ws.Range("E18").Value = "018-1234567"
$ws.Range("A19").Value = "Mei Fong"
$ws.Range("B19").Value = "kEkXMZDVFtFqFR"
$ws.Range("E19").Value = "019-2345678"
$ws.Range("A20").Value = "Raju Kumar"
$ws.Range("B20").Value = "q8GSAAPnlKe5z5"
$ws.Range("E20").Value = "010-3456789"
$ws.Range("A21").Value = "Zainal Abidin"
$ws.Range("B21").Value = "J3ksRMS8g22YGo"
$ws.Range("E21").Value = "011-4567890"
$ws.Range("A22").Value = "Amira binti Latif"
$ws.Range("B22").Value = "s5uro6raKM9eUh"
$ws.Range("E22").Value = "012-5678901"
$ws.Range("A23").Value = "Chia Wei Han"
$ws.Range("B23").Value = "wS4QiOIScAspiY"
$ws.Range("E23").Value = "013-6789012"
$ws.Range("A24").Value = "Suresh Kumar"
$ws.Range("B24").Value = "1fI0EladmVFdaO"
$ws.Range("E24").Value = "014-7890123"
$ws.Range("A25").Value = "Norhidayah"
$ws.Range("B25").Value = "2SbIreHqQRIlqm"
$ws.Range("E25").Value = "015-8901234"
$ws.Range("A26").Value = "Lim Wei Hui"
$ws.Range("B26").Value = "OebkWAaJpKgyU8"
$ws.Range("E26").Value = "016-9012345"
$ws.Range("A27").Value = "Arvind Raj"
$ws.Range("B27").Value = "ntEnbjY3tWLSLU"
$ws.Range("E27").Value = "017-0123456"
$ws.Range("A28").Value = "Safia binti Noor"
$ws.Range("B28").Value = "nZ9T1Bgcxj9eUw"
$ws.Range("E28").Value = "018-1234567"
$ws.Range("A29").Value = "Wong Kin Fai"
$ws.Range("B29").Value = "rXWaaeI2vOGrPu"
$ws.Range("E29").Value = "019-2345678"
$ws.Range("A30").Value = "Aina binti Osman"
$ws.Range("B30").Value = "NSvTwCID4y1RXQ"
$ws.Range("E30").Value = "010-3456789"
$ws.Range("A31").Value = "Kumar Narayanan"
$ws.Range("B31").Value = "ZLv2lnnPXBjjqb"
$ws.Range("E31").Value = "011-4567890"

Write-Output "Updated rows 2-31 for columns A, B, E"